# Salesforce test-cases workbook update:
# "Screenshot is saved in same folder as the report and working fine"
#
# The TC003 (Delete a Opportunity) and TC004 (Edit a Opportunity) rows both
# had a step that located the Opportunity record created earlier by its
# original name ("Try1"). Since the opportunity gets renamed during the
# edit test (TC004 renames it to "Tryedit") before it is used/deleted, the
# lookup step text is updated so it looks up the record by its current
# value ("Tryedit") instead of the stale "Try1" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStepText = 'click Opportunity record with value="Tryedit"'

# Row 4 -> TC003 "Delete a Opportunity" test case, step column J
$ws.Range("J4").Value = $newStepText

# Row 5 -> TC004 "Edit a Opportunity" test case, step column J
$ws.Range("J5").Value = $newStepText

# Row heights settled to their new auto-sized values after the edit
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 135
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 105

# Leave the saved selection on the last-edited cell, the way Excel would
# when the file was saved right after this edit.
[void]$ws.Range("T4").Select()
